$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Consolidate the report: the "Absent" column (H) should reflect
# whether the student was marked "Real" (present) that day (column E).
# If Real = 0 the student is Absent (H = 1), otherwise H = 0.
for ($r = 3; $r -le 21; $r++) {
    $real = $ws.Cells.Item($r, 5).Value()
    if ($real -eq 0) {
        $ws.Cells.Item($r, 8).Value = 1
    } else {
        $ws.Cells.Item($r, 8).Value = 0
    }
}
